# Trend epi - changed to only mon and tues; Web epi - added function to select certain columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Web epi: add new "adjust_func" value (keep_only_web_epi_cols) to the two web_epi rows ---
$ws.Range("C2").Value = "keep_only_web_epi_cols"
$ws.Range("C3").Value = "keep_only_web_epi_cols"

# --- Trend epi: restrict the day-of-week flags on rows 18-19 down to Mon (I) and Thu (L) only ---
$ws.Range("H18").ClearContents()
$ws.Range("J18").ClearContents()
$ws.Range("K18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H19").ClearContents()
$ws.Range("J19").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

# --- Update the saved cursor/selection position ---
[void]$ws.Range("C4").Select()
